$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 960.2
$ws.Range("I28").Value = 908.63635
$ws.Range("J28").Value = 1102
$ws.Range("K28").Value = 908.63635
$ws.Range("L28").Value = 1102
$ws.Range("M28").Value = -423.63635
$ws.Range("N28").Value = -2072

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1600
$ws.Range("J40").Value = 1600
$ws.Range("L40").Value = 1600
$ws.Range("N40").Value = -1950

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 325.53845
$ws.Range("I42").Value = 504.66666
$ws.Range("J42").Value = 172
$ws.Range("K42").Value = 1513.99998
$ws.Range("L42").Value = 516
$ws.Range("M42").Value = -1283.99998
$ws.Range("N42").Value = -976

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3967
$ws.Range("I43").Value = 2980
$ws.Range("J43").Value = 4296
$ws.Range("K43").Value = 2980
$ws.Range("L43").Value = 4296
$ws.Range("M43").Value = -2911
$ws.Range("N43").Value = -4434

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 23840.066
$ws.Range("J93").Value = 23840.066
$ws.Range("L93").Value = 23840.066
$ws.Range("N93").Value = -28832.066

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 9000
$ws.Range("I107").Value = 9000
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 9000
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -7080
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 872.30615
$ws.Range("J129").Value = 889.66315
$ws.Range("L129").Value = 2668.98945
$ws.Range("N129").Value = -12668.98945

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 3813.6667
$ws.Range("I131").Value = 2952.5
$ws.Range("J131").Value = 4797.857
$ws.Range("K131").Value = 8857.5
$ws.Range("L131").Value = 14393.571
$ws.Range("M131").Value = -3817.5
$ws.Range("N131").Value = -24473.571

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 32263644
$ws.Range("I132").Value = 38466844
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 115400532
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -115398002
$ws.Range("N132").Value = -26060

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 224075.22
$ws.Range("I141").Value = 251484.62
$ws.Range("J141").Value = 4800
$ws.Range("K141").Value = 754453.86
$ws.Range("L141").Value = 14400
$ws.Range("M141").Value = -749273.86
$ws.Range("N141").Value = -24760

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4741.125
$ws.Range("I32").Value = 4748.7383
$ws.Range("J32").Value = 4687.8335
$ws.Range("K32").Value = 4748.7383
$ws.Range("L32").Value = 4687.8335
$ws.Range("M32").Value = -4461.7383
$ws.Range("N32").Value = -5261.8335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1739.7693
$ws.Range("I97").Value = 1055.091
$ws.Range("K97").Value = 1055.091
$ws.Range("M97").Value = -559.0909999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H103").Value = 34404.363
$ws.Range("J103").Value = 34404.363
$ws.Range("L103").Value = 34404.363
$ws.Range("N103").Value = -36748.363

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2367.4583
$ws.Range("I132").Value = 1060.4445
$ws.Range("K132").Value = 3181.3335
$ws.Range("M132").Value = -651.3335000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H137").Value = 40068.332
$ws.Range("J137").Value = 40068.332
$ws.Range("L137").Value = 40068.332
$ws.Range("N137").Value = -50268.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 40333.332
$ws.Range("J59").Value = 40333.332
$ws.Range("L59").Value = 40333.332
$ws.Range("N59").Value = -42027.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 35068.184
$ws.Range("J95").Value = 35068.184
$ws.Range("L95").Value = 35068.184
$ws.Range("N95").Value = -40560.184

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H137").Value = 43480
$ws.Range("J137").Value = 43480
$ws.Range("L137").Value = 43480
$ws.Range("N137").Value = -53680

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 37037676
$ws.Range("I16").Value = 37037676
$ws.Range("K16").Value = 37037676
$ws.Range("M16").Value = -37037389

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2734.639
$ws.Range("J31").Value = 4484.8887
$ws.Range("L31").Value = 4484.8887
$ws.Range("N31").Value = -5074.8887

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2734.639
$ws.Range("J34").Value = 4484.8887
$ws.Range("L34").Value = 4484.8887
$ws.Range("N34").Value = -4888.8887

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 37037676
$ws.Range("I113").Value = 37037676
$ws.Range("K113").Value = 37037676
$ws.Range("M113").Value = -37035506

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 668622.4
$ws.Range("I5").Value = 626.5
$ws.Range("J5").Value = 954906.3
$ws.Range("K5").Value = 1879.5
$ws.Range("L5").Value = 2864718.9
$ws.Range("M5").Value = -1767.5
$ws.Range("N5").Value = -2864942.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 758.8200000000001
$ws.Range("I131").Value = 329.77777
$ws.Range("J131").Value = 801.25275
$ws.Range("K131").Value = 989.33331
$ws.Range("L131").Value = 2403.75825
$ws.Range("M131").Value = 4050.66669
$ws.Range("N131").Value = -12483.75825

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2430.5293
$ws.Range("I132").Value = 962.5
$ws.Range("J132").Value = 3735.4443
$ws.Range("K132").Value = 8662.5
$ws.Range("L132").Value = 33618.9987
$ws.Range("M132").Value = -6132.5
$ws.Range("N132").Value = -38678.9987

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 668622.4
$ws.Range("I135").Value = 626.5
$ws.Range("J135").Value = 954906.3
$ws.Range("K135").Value = 5638.5
$ws.Range("L135").Value = 8594156.700000001
$ws.Range("M135").Value = -3103.5
$ws.Range("N135").Value = -8599226.700000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 1556
$ws.Range("I137").Value = 1445
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 4335
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = 765
$ws.Range("N137").Value = -16200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1775.7646
$ws.Range("I102").Value = 1237.138
$ws.Range("K102").Value = 1237.138
$ws.Range("M102").Value = 384.8620000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2776.6875
$ws.Range("I132").Value = 1660.8948
$ws.Range("J132").Value = 4407.4614
$ws.Range("K132").Value = 4982.6844
$ws.Range("L132").Value = 13222.3842
$ws.Range("M132").Value = -2452.6844
$ws.Range("N132").Value = -18282.3842

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 37202
$ws.Range("J137").Value = 37202
$ws.Range("L137").Value = 37202
$ws.Range("N137").Value = -47402

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 708
$ws.Range("I35").Value = 708
$ws.Range("K35").Value = 708
$ws.Range("M35").Value = -372

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2274.6296
$ws.Range("I46").Value = 1791.75
$ws.Range("J46").Value = 2660.9333
$ws.Range("K46").Value = 1791.75
$ws.Range("L46").Value = 2660.9333
$ws.Range("M46").Value = -1603.75
$ws.Range("N46").Value = -3036.9333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 658.78
$ws.Range("I68").Value = 658.78
$ws.Range("K68").Value = 658.78
$ws.Range("M68").Value = 90.22000000000003

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 658.78
$ws.Range("I71").Value = 658.78
$ws.Range("K71").Value = 3293.9
$ws.Range("M71").Value = 450.1000000000004

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H124").Value = 40000
$ws.Range("J124").Value = 40000
$ws.Range("L124").Value = 40000
$ws.Range("N124").Value = -49820

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5196.346
$ws.Range("I132").Value = 1153.2727
$ws.Range("J132").Value = 8161.2666
$ws.Range("K132").Value = 3459.8181
$ws.Range("L132").Value = 24483.7998
$ws.Range("M132").Value = -929.8181
$ws.Range("N132").Value = -29543.7998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6947009
$ws.Range("I132").Value = 1776.5416
$ws.Range("J132").Value = 13892242
$ws.Range("K132").Value = 5329.6248
$ws.Range("L132").Value = 41676726
$ws.Range("M132").Value = -2799.6248
$ws.Range("N132").Value = -41681786
